{"js": "// Apply the LOT2023 syllabus update described by the diff.\n// Each edit is performed via a literal search + full-text replace so\n// existing run formatting (bold/italic) on the target run is preserved.\n\nasync function replaceOnce(body, oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText.substring(0, 60));\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1. Activation date bump.\nawait replaceOnce(\n  body,\n  \"Ativa\u00e7\u00e3o: 01/01/2018\",\n  \"Ativa\u00e7\u00e3o: 01/01/2025\"\n);\n\n// 2. Objetivos (PT) paragraph rewrite.\nawait replaceOnce(\n  body,\n  \"Demonstrar as principais etapas no desenvolvimento dos processos bioqu\u00edmicos industriais abordando aspectos bioqu\u00edmicos importantes na produ\u00e7\u00e3o de alimentos, e importantes metab\u00f3litos. Apresentar aos alunos uma vis\u00e3o das aplica\u00e7\u00f5es potenciais e estrat\u00e9gicas da biotecnologia moderna, incluindo as biorefinarias de lignocelul\u00f3sicos.Aprimorar o racioc\u00ednio e despertar o esp\u00edrito cr\u00edtico e a criatividade dos alunos na resolu\u00e7\u00e3o de problemas industriais envolvendo processos bioqu\u00edmicos.\",\n  \"Demonstrar as principais etapas no desenvolvimento dos processos bioqu\u00edmicos industriais abordando aspectos bioqu\u00edmicos importantes na produ\u00e7\u00e3o de alimentos, e importantes metab\u00f3litos. Apresentar aos alunos uma vis\u00e3o das aplica\u00e7\u00f5es potenciais e estrat\u00e9gicas da biotecnologia moderna, incluindo aspectos bioqu\u00edmicos de bioprocessos envolvendo a utiliza\u00e7\u00e3o de hidrolisados lignocelul\u00f3sicos e suas aplica\u00e7\u00f5es tecnol\u00f3gicas. Aprimorar o racioc\u00ednio e despertar o esp\u00edrito cr\u00edtico e a criatividade dos alunos na resolu\u00e7\u00e3o de problemas industriais envolvendo processos bioqu\u00edmicos.\"\n);\n\n// 3. Objetivos (EN) paragraph: previously an empty italic run, now filled in.\n{\n  const paras = body.paragraphs;\n  paras.load(\"text\");\n  await context.sync();\n  let target = null;\n  for (let i = 0; i < paras.items.length; i++) {\n    if (paras.items[i].text === \"\") {\n      // The empty italic paragraph directly follows the Objetivos (PT) paragraph.\n      const prevText = i > 0 ? paras.items[i - 1].text : \"\";\n      if (prevText.indexOf(\"Aprimorar o racioc\u00ednio\") !== -1) {\n        target = paras.items[i];\n        break;\n      }\n    }\n  }\n  if (!target) {\n    throw new Error(\"Could not locate empty Objetivos (EN) paragraph\");\n  }\n  target.insertText(\n    \"Demonstrate the main steps in the development of industrial biochemical processes, addressing important biochemical aspects in food production, and important metabolites. Present students with a vision of the potential and strategic applications of modern biotechnology, including biochemical aspects of bioprocesses involving the use of lignocellulosic hydrolysates and their technological applications. Improve reasoning and awaken students' critical spirit and creativity when solving industrial problems involving biochemical processes\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 4. Programa resumido (PT).\nawait replaceOnce(\n  body,\n  \"Introdu\u00e7\u00e3o aos processos bioqu\u00edmicos industriais que incluem o processamento de alimentos, e importantes metab\u00f3litos, a manufatura de soros e vacinas, e os conceitos modernos de bioenergia e biorrefinarias.\",\n  \"Introdu\u00e7\u00e3o aos processos bioqu\u00edmicos industriais que incluem o processamento de alimentos, e importantes metab\u00f3litos, a manufatura de bioprodutos, e os aspectos bioqu\u00edmicos de bioprocessos envolvendo bioenergia e biorrefinarias.\"\n);\n\n// 5. Programa resumido (EN).\nawait replaceOnce(\n  body,\n  \"Introduction to industrial biochemical processes which include the food processing and important metabolites, the manufacture of serum and vaccines, the modern concepts of bioenergy and biorefineries.\",\n  \"Introduction to industrial biochemical processes that include food processing and important metabolites, the manufacture of bioproducts, and the biochemical aspects of bioprocesses involving bioenergy and biorefineries.\"\n);\n\n// 6. Programa (PT).\nawait replaceOnce(\n  body,\n  \"1. Introdu\u00e7\u00e3o ao Processamento de alimentos: tipos de ind\u00fastria de alimentos, mat\u00e9riasprimas,fases doprocessamento de produtos aliment\u00edcios, conserva\u00e7\u00e3o/altera\u00e7\u00f5es de alimentos, microbiologia dealimentos, altera\u00e7\u00f5es bioqu\u00edmicas em alimentos (oxida\u00e7\u00e3o de lip\u00eddeos, antioxidantes, escurecimentoenzim\u00e1tico e n\u00e3o enzim\u00e1tico), aflatoxinas, conservantes qu\u00edmicos, toxicantes naturais.2. Discuss\u00e3o e apresenta\u00e7\u00e3o sobre aspectos bioqu\u00edmicos importantes na produ\u00e7\u00e3o de metab\u00f3litos por microrganismos e  estudo de casos .3. Manufatura de soros e vacinas M\u00e9todosindustriais para a produ\u00e7\u00e3o de soros e vacinas 4.Biotecnologia de lignocelul\u00f3sicos: Separa\u00e7\u00e3o e fermenta\u00e7\u00e3o das fra\u00e7\u00f5es e principais processosbioqu\u00edmicos envolvendo materiais lignocelul\u00f3sicos.5. Bioenergia, biocombust\u00edveis e biorrefinarias.\",\n  \"1. Introdu\u00e7\u00e3o: abordagem geral dos princ\u00edpios bioqu\u00edmicos aplicados em diferentes processos2. Processos bioqu\u00edmicos aplicados \u00e0 ind\u00fastria de alimentos: tipos de ind\u00fastria de alimentos, mat\u00e9ria primas, fases do processamento de produtos aliment\u00edcios, conserva\u00e7\u00e3o/altera\u00e7\u00f5es de alimentos 3. Processos bioqu\u00edmicos nas ind\u00fastrias de processamento de produtos l\u00e1cteos, frutas e hortali\u00e7as, cacau, produtos gordurosos, produtos desidratados 4. Discuss\u00e3o e apresenta\u00e7\u00e3o sobre aspectos bioqu\u00edmicos na produ\u00e7\u00e3o de bioprodutos: solventes, \u00e1cidos org\u00e2nicos, f\u00e1rmacos, soros e vacinas, bioinsumos agr\u00edcolas e outros de import\u00e2ncia industrial5. Bioenergia e biorrefinarias: aspectos bioqu\u00edmicos de bioprocessos envolvendo a utiliza\u00e7\u00e3o de hidrolisados lignocelul\u00f3sicos e suas aplica\u00e7\u00f5es tecnol\u00f3gicas\"\n);\n\n// 7. Programa (EN).\nawait replaceOnce(\n  body,\n  \"1. Food processing: kinds of food industries, rawmaterials,processing phases of food products,conservation/changes of foods, microbiology of foods,biochemical changes in foods (oxidation oflipids, antioxidants, enzymatic and nonenzymaticdarkening), aflatoxins, chemical preservatives,natural toxicants.2. Discussion and presentation on important biochemical aspects in the production of metabolites by microorganisms and case studies.3. Production of serum and vaccines: industrial methods for production of serum and vaccines \u20134. Biotechnology of lignocellulose: separation and fermentation of fraction and main biochemicalprocesses involving lignocellulosic materials.5. Bioenergy, biofuels and biorefineries\",\n  \"1. Introduction: general approach to biochemical principles applied in different processes2. Biochemical processes applied to the food industry: types of food industry, raw materials, stages of food processing, food conservation/modifications3. Biochemical processes in the processing industries of dairy products, fruits and vegetables, cocoa, fatty products, dehydrated products4. Discussion and presentation on biochemical aspects in the production of bioproducts: solvents, organic acids, pharmaceuticals, serums and vaccines, agricultural bioinputs and others of industrial importance5. Bioenergy and biorefineries: biochemical aspects of bioprocesses involving the use of lignocellulosic hydrolysates and their technological applications\"\n);\n\n// 8. Norma de recupera\u00e7\u00e3o: missing space fix.\nawait replaceOnce(\n  body,\n  \"A recupera\u00e7\u00e3o ser\u00e1 feita por meio de uma prova escrita (PR) e a m\u00e9dia de recupera\u00e7\u00e3o (MR)ser\u00e1 calculada como MR=(NF+PR)/2\",\n  \"A recupera\u00e7\u00e3o ser\u00e1 feita por meio de uma prova escrita (PR) e a m\u00e9dia de recupera\u00e7\u00e3o (MR) ser\u00e1 calculada como MR=(NF+PR)/2\"\n);\n\n// 9. Bibliografia: collapse the 7 references (joined by line breaks) into a\n// single paragraph of new reference text.\n{\n  const paras = body.paragraphs;\n  paras.load(\"text\");\n  await context.sync();\n  let target = null;\n  for (let i = 0; i < paras.items.length; i++) {\n    if (paras.items[i].text.indexOf(\"BORZANI\") !== -1) {\n      target = paras.items[i];\n      break;\n    }\n  }\n  if (!target) {\n    throw new Error(\"Could not locate Bibliografia paragraph\");\n  }\n  target.insertText(\n    \"GAVA, A. J.; SILVA, C. A. B.; FRIAS, J. R. B. Tecnologia de alimentos - princ\u00edpios e aplica\u00e7\u00f5es. S\u00e3o Paulo, Nobel, 2008. ISBN-13: 9788521313823.LIMA, U. A. Biotecnologia Industrial: Processos Fermentativos e Enzim\u00e1ticos. Volume 3. S\u00e3o Paulo: Editora Edgard Bl\u00fccher Ltda, 2019. ISBN 9788521214571.Moraes, I. O. Biotecnologia Industrial: Biotecnologia na produ\u00e7\u00e3o de alimentos. Volume 4. S\u00e3o Paulo: Editora Edgard Bl\u00fccher Ltda, 2021. ISBN 9786555061529.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Apply the LOT2023 syllabus update described by the diff using Word COM\n# interop (Find/Replace for plain phrase swaps, Paragraphs.Item(...).Range.Text\n# for the two spots that need a full-paragraph rewrite: the previously-empty\n# English \"Objetivos\" paragraph and the Bibliografia paragraph whose seven\n# line-broken references collapse into one new reference list).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n# 1. Activation date bump.\nReplace-Text \"Ativa\u00e7\u00e3o: 01/01/2018\" \"Ativa\u00e7\u00e3o: 01/01/2025\"\n\n# 2. Objetivos (PT) paragraph rewrite.\nReplace-Text `\n    \"Demonstrar as principais etapas no desenvolvimento dos processos bioqu\u00edmicos industriais abordando aspectos bioqu\u00edmicos importantes na produ\u00e7\u00e3o de alimentos, e importantes metab\u00f3litos. Apresentar aos alunos uma vis\u00e3o das aplica\u00e7\u00f5es potenciais e estrat\u00e9gicas da biotecnologia moderna, incluindo as biorefinarias de lignocelul\u00f3sicos.Aprimorar o racioc\u00ednio e despertar o esp\u00edrito cr\u00edtico e a criatividade dos alunos na resolu\u00e7\u00e3o de problemas industriais envolvendo processos bioqu\u00edmicos.\" `\n    \"Demonstrar as principais etapas no desenvolvimento dos processos bioqu\u00edmicos industriais abordando aspectos bioqu\u00edmicos importantes na produ\u00e7\u00e3o de alimentos, e importantes metab\u00f3litos. Apresentar aos alunos uma vis\u00e3o das aplica\u00e7\u00f5es potenciais e estrat\u00e9gicas da biotecnologia moderna, incluindo aspectos bioqu\u00edmicos de bioprocessos envolvendo a utiliza\u00e7\u00e3o de hidrolisados lignocelul\u00f3sicos e suas aplica\u00e7\u00f5es tecnol\u00f3gicas. Aprimorar o racioc\u00ednio e despertar o esp\u00edrito cr\u00edtico e a criatividade dos alunos na resolu\u00e7\u00e3o de problemas industriais envolvendo processos bioqu\u00edmicos.\"\n\n# 3. Objetivos (EN) paragraph: previously an empty italic run, now filled in.\n#    Locate it as the (only) essentially-empty paragraph carrying italic\n#    formatting, which survives renumbering regardless of the earlier edits.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Length -le 1 -and $p.Range.Font.Italic) {\n        $p.Range.Text = \"Demonstrate the main steps in the development of industrial biochemical processes, addressing important biochemical aspects in food production, and important metabolites. Present students with a vision of the potential and strategic applications of modern biotechnology, including biochemical aspects of bioprocesses involving the use of lignocellulosic hydrolysates and their technological applications. Improve reasoning and awaken students' critical spirit and creativity when solving industrial problems involving biochemical processes\"\n        break\n    }\n}\n\n# 4. Programa resumido (PT).\nReplace-Text `\n    \"Introdu\u00e7\u00e3o aos processos bioqu\u00edmicos industriais que incluem o processamento de alimentos, e importantes metab\u00f3litos, a manufatura de soros e vacinas, e os conceitos modernos de bioenergia e biorrefinarias.\" `\n    \"Introdu\u00e7\u00e3o aos processos bioqu\u00edmicos industriais que incluem o processamento de alimentos, e importantes metab\u00f3litos, a manufatura de bioprodutos, e os aspectos bioqu\u00edmicos de bioprocessos envolvendo bioenergia e biorrefinarias.\"\n\n# 5. Programa resumido (EN).\nReplace-Text `\n    \"Introduction to industrial biochemical processes which include the food processing and important metabolites, the manufacture of serum and vaccines, the modern concepts of bioenergy and biorefineries.\" `\n    \"Introduction to industrial biochemical processes that include food processing and important metabolites, the manufacture of bioproducts, and the biochemical aspects of bioprocesses involving bioenergy and biorefineries.\"\n\n# 6. Programa (PT).\nReplace-Text `\n    \"1. Introdu\u00e7\u00e3o ao Processamento de alimentos: tipos de ind\u00fastria de alimentos, mat\u00e9riasprimas,fases doprocessamento de produtos aliment\u00edcios, conserva\u00e7\u00e3o/altera\u00e7\u00f5es de alimentos, microbiologia dealimentos, altera\u00e7\u00f5es bioqu\u00edmicas em alimentos (oxida\u00e7\u00e3o de lip\u00eddeos, antioxidantes, escurecimentoenzim\u00e1tico e n\u00e3o enzim\u00e1tico), aflatoxinas, conservantes qu\u00edmicos, toxicantes naturais.2. Discuss\u00e3o e apresenta\u00e7\u00e3o sobre aspectos bioqu\u00edmicos importantes na produ\u00e7\u00e3o de metab\u00f3litos por microrganismos e  estudo de casos .3. Manufatura de soros e vacinas M\u00e9todosindustriais para a produ\u00e7\u00e3o de soros e vacinas 4.Biotecnologia de lignocelul\u00f3sicos: Separa\u00e7\u00e3o e fermenta\u00e7\u00e3o das fra\u00e7\u00f5es e principais processosbioqu\u00edmicos envolvendo materiais lignocelul\u00f3sicos.5. Bioenergia, biocombust\u00edveis e biorrefinarias.\" `\n    \"1. Introdu\u00e7\u00e3o: abordagem geral dos princ\u00edpios bioqu\u00edmicos aplicados em diferentes processos2. Processos bioqu\u00edmicos aplicados \u00e0 ind\u00fastria de alimentos: tipos de ind\u00fastria de alimentos, mat\u00e9ria primas, fases do processamento de produtos aliment\u00edcios, conserva\u00e7\u00e3o/altera\u00e7\u00f5es de alimentos 3. Processos bioqu\u00edmicos nas ind\u00fastrias de processamento de produtos l\u00e1cteos, frutas e hortali\u00e7as, cacau, produtos gordurosos, produtos desidratados 4. Discuss\u00e3o e apresenta\u00e7\u00e3o sobre aspectos bioqu\u00edmicos na produ\u00e7\u00e3o de bioprodutos: solventes, \u00e1cidos org\u00e2nicos, f\u00e1rmacos, soros e vacinas, bioinsumos agr\u00edcolas e outros de import\u00e2ncia industrial5. Bioenergia e biorrefinarias: aspectos bioqu\u00edmicos de bioprocessos envolvendo a utiliza\u00e7\u00e3o de hidrolisados lignocelul\u00f3sicos e suas aplica\u00e7\u00f5es tecnol\u00f3gicas\"\n\n# 7. Programa (EN).\nReplace-Text `\n    \"1. Food processing: kinds of food industries, rawmaterials,processing phases of food products,conservation/changes of foods, microbiology of foods,biochemical changes in foods (oxidation oflipids, antioxidants, enzymatic and nonenzymaticdarkening), aflatoxins, chemical preservatives,natural toxicants.2. Discussion and presentation on important biochemical aspects in the production of metabolites by microorganisms and case studies.3. Production of serum and vaccines: industrial methods for production of serum and vaccines \u20134. Biotechnology of lignocellulose: separation and fermentation of fraction and main biochemicalprocesses involving lignocellulosic materials.5. Bioenergy, biofuels and biorefineries\" `\n    \"1. Introduction: general approach to biochemical principles applied in different processes2. Biochemical processes applied to the food industry: types of food industry, raw materials, stages of food processing, food conservation/modifications3. Biochemical processes in the processing industries of dairy products, fruits and vegetables, cocoa, fatty products, dehydrated products4. Discussion and presentation on biochemical aspects in the production of bioproducts: solvents, organic acids, pharmaceuticals, serums and vaccines, agricultural bioinputs and others of industrial importance5. Bioenergy and biorefineries: biochemical aspects of bioprocesses involving the use of lignocellulosic hydrolysates and their technological applications\"\n\n# 8. Norma de recupera\u00e7\u00e3o: missing space fix.\nReplace-Text `\n    \"A recupera\u00e7\u00e3o ser\u00e1 feita por meio de uma prova escrita (PR) e a m\u00e9dia de recupera\u00e7\u00e3o (MR)ser\u00e1 calculada como MR=(NF+PR)/2\" `\n    \"A recupera\u00e7\u00e3o ser\u00e1 feita por meio de uma prova escrita (PR) e a m\u00e9dia de recupera\u00e7\u00e3o (MR) ser\u00e1 calculada como MR=(NF+PR)/2\"\n\n# 9. Bibliografia: collapse the 7 references (joined by line breaks) into a\n#    single paragraph of new reference text.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*BORZANI*\") {\n        $p.Range.Text = \"GAVA, A. J.; SILVA, C. A. B.; FRIAS, J. R. B. Tecnologia de alimentos - princ\u00edpios e aplica\u00e7\u00f5es. S\u00e3o Paulo, Nobel, 2008. ISBN-13: 9788521313823.LIMA, U. A. Biotecnologia Industrial: Processos Fermentativos e Enzim\u00e1ticos. Volume 3. S\u00e3o Paulo: Editora Edgard Bl\u00fccher Ltda, 2019. ISBN 9788521214571.Moraes, I. O. Biotecnologia Industrial: Biotecnologia na produ\u00e7\u00e3o de alimentos. Volume 4. S\u00e3o Paulo: Editora Edgard Bl\u00fccher Ltda, 2021. ISBN 9786555061529.\"\n        break\n    }\n}\n"}
